# The doc ends with a paragraph reading "Here is the link to my GitHub".
# Append a blank paragraph and then a paragraph containing the repo URL,
# both inheriting that paragraph's formatting (Times New Roman / 12pt /
# first-line indent), matching how Word extends a run of similarly
# formatted paragraphs when you press Enter at the end of the document.

$d = $word.ActiveDocument

$anchor = $d.Paragraphs.Last
if ($anchor.Range.Text -notmatch "Here is the link to my GitHub") {
    # Fall back to locating it explicitly if it is not the last paragraph.
    $find = $d.Content
    $find.Find.Execute("Here is the link to my GitHub", $true, $false, $false,
                        $false, $false, $true, 1, $false, "", 0)
    $anchor = $find.Paragraphs.Last
}

# Blank paragraph right after the GitHub intro line.
$anchor.Range.InsertParagraphAfter()

# Paragraph holding the repository URL.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "https://github.com/DHadi-5722/csd-310.git"
